$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet "Sheet2" right after "Sheet1" (it becomes the active tab)
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Column widths for the new sheet
$ws2.Columns.Item(1).ColumnWidth = 30.998697916666668
$ws2.Columns.Item(2).ColumnWidth = 79.99869791666667
$ws2.Columns.Item(3).ColumnWidth = 7.330729166666667

# Header row (re-uses the same strings already used on Sheet1)
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Homework 1"
$ws2.Range("C1").Value = "Grade"
$ws2.Range("D1").Value = "Homework 2"
$ws2.Range("E1").Value = "Grade"
$ws2.Range("F1").Value = "Homework 3"
$ws2.Range("G1").Value = "Grade"

# Grading data for Homework 1 of this workshop (entered in the order the
# unique names/comments first show up so the shared string table matches)
$ws2.Range("A2").Value = "Chanbormey Ouk"
$ws2.Range("A3").Value = "Clara Zuccarino"
$ws2.Range("A4").Value = "Emma Guitera"
$ws2.Range("B4").Value = "Good effort but some answers were incorrect…"
$ws2.Range("A5").Value = "Esteban Filippa-Ocampo"
$ws2.Range("B2").Value = "Good job (only unable to do bonus question). Next time please comment answers"
$ws2.Range("B3").Value = "Good job (only unable to do bonus question)"
$ws2.Range("A6").Value = "Ishankaur Khalsa"
$ws2.Range("A7").Value = "Jane Wilkinson"
$ws2.Range("B7").Value = "Good effort, some answers were incorrect but attempted the bonus question"
$ws2.Range("A8").Value = "Kevin Alejandro Gomez Anguiano"
$ws2.Range("A9").Value = "Luciana Markstein"
$ws2.Range("B9").Value = "Great job, even did bonus question!"
$ws2.Range("A10").Value = "Paniz Jamali"
$ws2.Range("B10").Value = "Sent us a PHG file…"
$ws2.Range("A11").Value = "Samia Firmino Pinto"
$ws2.Range("B11").Value = "Good effort, not the most straighforard way to do some things but most answers were correct!"
$ws2.Range("A12").Value = "Sebastian Zünd"
$ws2.Range("A13").Value = "Willow Sylvester"
$ws2.Range("A14").Value = "Xinyu Liu"
$ws2.Range("A15").Value = "Yixin Du"
$ws2.Range("A16").Value = "Nameetta Nierakkal"

# Remaining (duplicate) comments
$ws2.Range("B5").Value = "Good job (only unable to do bonus question)"
$ws2.Range("B6").Value = "Good job (only unable to do bonus question). Next time please comment answers"
$ws2.Range("B8").Value = "Good job (only unable to do bonus question). Next time please comment answers"
$ws2.Range("B12").Value = "Good job (only unable to do bonus question)"
$ws2.Range("B13").Value = "Great job, even did bonus question!"
$ws2.Range("B14").Value = "Good effort, some answers were incorrect but attempted the bonus question"
$ws2.Range("B15").Value = "Good effort, some answers were incorrect but attempted the bonus question"
$ws2.Range("B16").Value = "Good effort but some answers were incorrect…"

# Grades
$ws2.Range("C2").Value = 6
$ws2.Range("C3").Value = 6
$ws2.Range("C4").Value = 5.5
$ws2.Range("C5").Value = 6
$ws2.Range("C6").Value = 6
$ws2.Range("C7").Value = 5.75
$ws2.Range("C8").Value = 6
$ws2.Range("C9").Value = 6
$ws2.Range("C11").Value = 5.5
$ws2.Range("C12").Value = 6
$ws2.Range("C13").Value = 6
$ws2.Range("C14").Value = 5.75
$ws2.Range("C15").Value = 5.75
$ws2.Range("C16").Value = 5.5

# Update view state: Sheet1 is no longer the selected tab, and its selection
# becomes the whole of row 1 (as if the row header had been clicked); Sheet2
# becomes the active tab with cell B17 selected.
$ws1.Rows.Item(1).Select()
$ws2.Range("B17").Select()
